# Auto update Excel log
# Appends newly-collected sensor readings to the PIR, Humidity, and
# Temperature worksheets of the SeniorConnect master log.

$wb = $excel.ActiveWorkbook

# Helper: write a string value to a cell while preventing Excel's
# automatic text-to-number/date conversion (e.g. "2026-01-28" becoming a
# date serial, or "88.4%" becoming a numeric percentage). We briefly mark
# the cell as Text, assign the literal string, then restore the cell's
# style to "Normal" so no stray formatting is left behind.
function Set-TextValue {
    param($Sheet, $Row, $Col, $Text)
    $cell = $Sheet.Cells.Item($Row, $Col)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# PIR sheet: append rows 107-119 (columns: Date, Timestamp, Hour, Location,
# Value, Status)
# ---------------------------------------------------------------------------
$wsPIR = $wb.Worksheets.Item("PIR")
$pirData = @(
        @("2026-01-28", "15:37:15", "15:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-28", "15:37:16", "15:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-28", "15:37:21", "15:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-28", "15:37:26", "15:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-28", "15:37:31", "15:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-28", "15:37:36", "15:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-28", "15:37:41", "15:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-28", "15:37:46", "15:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-28", "15:37:51", "15:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-28", "15:37:56", "15:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-28", "15:38:01", "15:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-28", "15:38:06", "15:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-28", "15:38:11", "15:00", "Bathroom", "No Motion", "Inactive")
)
$pirStartRow = 107
for ($i = 0; $i -lt $pirData.Count; $i++) {
    $row = $pirStartRow + $i
    $values = $pirData[$i]
    Set-TextValue $wsPIR $row 1 $values[0]
    $wsPIR.Cells.Item($row, 2).Value = $values[1]
    $wsPIR.Cells.Item($row, 3).Value = $values[2]
    $wsPIR.Cells.Item($row, 4).Value = $values[3]
    $wsPIR.Cells.Item($row, 5).Value = $values[4]
    $wsPIR.Cells.Item($row, 6).Value = $values[5]
}

# ---------------------------------------------------------------------------
# Humidity sheet: append rows 109-120 (columns: Date, Timestamp, Hour,
# Location, Value, Status). The Value column holds percentage text such as
# "88.4%" which must stay text rather than becoming a numeric percentage.
# ---------------------------------------------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")
$humidityData = @(
        @("2026-01-28", "15:37:15", "15:00", "Bathroom", "88.4%", "Active"),
        @("2026-01-28", "15:37:16", "15:00", "Bathroom", "88.4%", "Active"),
        @("2026-01-28", "15:37:19", "15:00", "Bathroom", "87.5%", "Active"),
        @("2026-01-28", "15:37:23", "15:00", "Bathroom", "88.4%", "Active"),
        @("2026-01-28", "15:37:31", "15:00", "Bathroom", "88.4%", "Active"),
        @("2026-01-28", "15:37:35", "15:00", "Bathroom", "88.4%", "Active"),
        @("2026-01-28", "15:37:43", "15:00", "Bathroom", "88.4%", "Active"),
        @("2026-01-28", "15:37:47", "15:00", "Bathroom", "87.5%", "Active"),
        @("2026-01-28", "15:37:51", "15:00", "Bathroom", "88.4%", "Active"),
        @("2026-01-28", "15:37:55", "15:00", "Bathroom", "88.4%", "Active"),
        @("2026-01-28", "15:38:03", "15:00", "Bathroom", "88.4%", "Active"),
        @("2026-01-28", "15:38:12", "15:00", "Bathroom", "88.4%", "Active")
)
$humidityStartRow = 109
for ($i = 0; $i -lt $humidityData.Count; $i++) {
    $row = $humidityStartRow + $i
    $values = $humidityData[$i]
    Set-TextValue $wsHumidity $row 1 $values[0]
    $wsHumidity.Cells.Item($row, 2).Value = $values[1]
    $wsHumidity.Cells.Item($row, 3).Value = $values[2]
    $wsHumidity.Cells.Item($row, 4).Value = $values[3]
    Set-TextValue $wsHumidity $row 5 $values[4]
    $wsHumidity.Cells.Item($row, 6).Value = $values[5]
}

# ---------------------------------------------------------------------------
# Temperature sheet: append rows 109-120 (columns: Date, Timestamp, Hour,
# Location, Value, Status)
# ---------------------------------------------------------------------------
$wsTemperature = $wb.Worksheets.Item("Temperature")
$temperatureData = @(
        @("2026-01-28", "15:37:15", "15:00", "Bathroom", "22.9C", "Active"),
        @("2026-01-28", "15:37:16", "15:00", "Bathroom", "22.9C", "Active"),
        @("2026-01-28", "15:37:20", "15:00", "Bathroom", "22.9C", "Active"),
        @("2026-01-28", "15:37:24", "15:00", "Bathroom", "22.9C", "Active"),
        @("2026-01-28", "15:37:32", "15:00", "Bathroom", "22.9C", "Active"),
        @("2026-01-28", "15:37:36", "15:00", "Bathroom", "22.9C", "Active"),
        @("2026-01-28", "15:37:44", "15:00", "Bathroom", "22.9C", "Active"),
        @("2026-01-28", "15:37:48", "15:00", "Bathroom", "22.9C", "Active"),
        @("2026-01-28", "15:37:52", "15:00", "Bathroom", "22.9C", "Active"),
        @("2026-01-28", "15:37:56", "15:00", "Bathroom", "22.9C", "Active"),
        @("2026-01-28", "15:38:04", "15:00", "Bathroom", "22.9C", "Active"),
        @("2026-01-28", "15:38:12", "15:00", "Bathroom", "22.9C", "Active")
)
$temperatureStartRow = 109
for ($i = 0; $i -lt $temperatureData.Count; $i++) {
    $row = $temperatureStartRow + $i
    $values = $temperatureData[$i]
    Set-TextValue $wsTemperature $row 1 $values[0]
    $wsTemperature.Cells.Item($row, 2).Value = $values[1]
    $wsTemperature.Cells.Item($row, 3).Value = $values[2]
    $wsTemperature.Cells.Item($row, 4).Value = $values[3]
    $wsTemperature.Cells.Item($row, 5).Value = $values[4]
    $wsTemperature.Cells.Item($row, 6).Value = $values[5]
}
